$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.545.79"
$ws.Cells.Item(2, 5).Value = "  +1.69%  "
$ws.Cells.Item(3, 4).Value = "2.544.58"
$ws.Cells.Item(3, 5).Value = "  +4.83%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).Value = "'568.91"
$ws.Cells.Item(5, 5).Value = "  +2.14%  "
$ws.Cells.Item(6, 4).Value = "'150.56"
$ws.Cells.Item(6, 5).Value = "  +8.43%  "
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(8, 4).Value = "'0.586"
$ws.Cells.Item(8, 5).Value = "  +0.10%  "
$ws.Cells.Item(9, 4).Value = "2.544.63"
$ws.Cells.Item(9, 5).Value = "  +4.90%  "
$ws.Cells.Item(10, 4).Value = "'0.106"
$ws.Cells.Item(10, 5).Value = "  +1.81%  "
$ws.Cells.Item(11, 4).Value = "'5.71"
$ws.Cells.Item(11, 5).Value = "  -1.23%  "
$ws.Cells.Item(12, 5).Value = "  +1.04%  "
$ws.Cells.Item(13, 4).Value = "'0.358"
$ws.Cells.Item(13, 5).Value = "  +2.82%  "
$ws.Cells.Item(14, 4).Value = "'28.25"
$ws.Cells.Item(14, 5).Value = "  +8.08%  "
$ws.Cells.Item(15, 4).Value = "3.002.47"
$ws.Cells.Item(15, 5).Value = "  +4.96%  "
$ws.Cells.Item(16, 4).Value = "63.458.36"
$ws.Cells.Item(16, 5).Value = "  +1.74%  "
$ws.Cells.Item(17, 5).Value = "  +2.00%  "
$ws.Cells.Item(18, 4).Value = "2.524.65"
$ws.Cells.Item(18, 5).Value = "  +3.85%  "
$ws.Cells.Item(19, 4).Value = "'11.61"
$ws.Cells.Item(19, 5).Value = "  +3.82%  "
$ws.Cells.Item(20, 4).Value = "'339.46"
$ws.Cells.Item(20, 5).Value = "  -1.77%  "
$ws.Cells.Item(21, 4).Value = "'4.35"
$ws.Cells.Item(21, 5).Value = "  +3.23%  "
$ws.Cells.Item(22, 4).Value = "'6.81"
$ws.Cells.Item(22, 5).Value = "  -0.59%  "
$ws.Cells.Item(23, 5).Value = "  +0.16%  "
$ws.Cells.Item(24, 5).Value = "  +1.31%  "
$ws.Cells.Item(25, 5).Value = "  -1.30%  "
$ws.Cells.Item(26, 4).Value = "'1.56"
$ws.Cells.Item(26, 5).Value = "  +15.43%  "
$ws.Cells.Item(27, 5).Value = "  +3.51%  "
$ws.Cells.Item(28, 4).Value = "'8.47"
$ws.Cells.Item(28, 5).Value = "  +4.03%  "
$ws.Cells.Item(29, 5).Value = "  -0.16%  "
$ws.Cells.Item(30, 5).Value = "  +10.62%  "
$ws.Cells.Item(31, 4).Value = "0.0₃0823"
$ws.Cells.Item(31, 5).Value = "  +4.70%  "
$ws.Cells.Item(32, 4).Value = "'1.88"
$ws.Cells.Item(32, 5).Value = "  +3.91%  "
$ws.Cells.Item(33, 5).Value = "  +3.24%  "
$ws.Cells.Item(34, 5).Value = "  +10.55%  "
$ws.Cells.Item(35, 4).Value = "'420.48"
$ws.Cells.Item(35, 5).Value = "  +10.72%  "
$ws.Cells.Item(36, 5).Value = "  +2.37%  "
$ws.Cells.Item(37, 4).Value = "'19.05"
$ws.Cells.Item(37, 5).Value = "  +2.51%  "
$ws.Cells.Item(38, 4).Value = "'4.43"
$ws.Cells.Item(38, 5).Value = "  -0.49%  "
$ws.Cells.Item(39, 5).Value = "  -0.01%  "
$ws.Cells.Item(40, 4).Value = "'1.78"
$ws.Cells.Item(40, 5).Value = "  +5.11%  "
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 5).Value = "  +0.04%  "
$ws.Cells.Item(42, 4).Value = "'39.54"
$ws.Cells.Item(42, 5).Value = "  +0.54%  "
$ws.Cells.Item(43, 4).Value = "'153.85"
$ws.Cells.Item(43, 5).Value = "  +6.10%  "
$ws.Cells.Item(44, 4).Value = "'3.78"
$ws.Cells.Item(44, 5).Value = "  +3.10%  "
$ws.Cells.Item(45, 4).Value = "'20.85"
$ws.Cells.Item(45, 5).Value = "  +0.19%  "
$ws.Cells.Item(46, 4).Value = "'0.611"
$ws.Cells.Item(46, 5).Value = "  +3.22%  "
$ws.Cells.Item(47, 4).Value = "'0.0528"
$ws.Cells.Item(47, 5).Value = "  +1.53%  "
$ws.Cells.Item(48, 4).Value = "'0.0966"
$ws.Cells.Item(48, 5).Value = "  +1.24%  "
$ws.Cells.Item(49, 4).Value = "'0.0239"
$ws.Cells.Item(49, 5).Value = "  +7.08%  "
$ws.Cells.Item(50, 4).Value = "'18.65"
$ws.Cells.Item(50, 5).Value = "  +4.40%  "
$ws.Cells.Item(51, 4).Value = "'1.81"
$ws.Cells.Item(51, 5).Value = "  +5.81%  "
